$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Piso 2": update rebar spacing data (columns F y G, filas 10 a 24)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Piso 2")

$ws3.Range("G10").Value = 400

$ws3.Range("F11").Value = 6
$ws3.Range("G11").Value = 400

$ws3.Range("G12").Value = 200

$ws3.Range("G13").Value = 200

$ws3.Range("F14").Value = "2a"
$ws3.Range("G14").Value = 300

$ws3.Range("F15").Value = "5b"
$ws3.Range("G15").Value = 300

$ws3.Range("F16").Value = "5b"
$ws3.Range("G16").Value = 300

$ws3.Range("F17").Value = "2a"
$ws3.Range("G17").Value = 300

$ws3.Range("F18").Value = "2a"
$ws3.Range("G18").Value = 300

$ws3.Range("F19").Value = "2a"
$ws3.Range("G19").Value = 300

$ws3.Range("F21").Value = "2a"
$ws3.Range("G21").Value = 300

$ws3.Range("F22").Value = "2a"
$ws3.Range("G22").Value = 300

$ws3.Range("F23").Value = "2a"
$ws3.Range("G23").Value = 300

$ws3.Range("F24").Value = "5b"

# ---------------------------------------------------------------------------
# Actualizar selecciones / vistas activas de cada hoja
# ---------------------------------------------------------------------------

# "Piso 1": mover la seleccion de B2:H22 a F15
$ws2 = $wb.Worksheets.Item("Piso 1")
$ws2.Activate()
$ws2.Range("F15").Select()

# "Piso 2": mover la seleccion de F10 a G24 (y ya no es la hoja activa)
$ws3.Activate()
$ws3.Range("G24").Select()

# "Piso -1": pasa a ser la hoja activa/seleccionada, viendo desde B1
$ws1 = $wb.Worksheets.Item("Piso -1")
$ws1.Activate()
$ws1.Range("B1").Select()
$ws1.Range("B32").Select()
